$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 26: update title (D26)
$ws.Range("D26").Value = "인공지능 음성 생성 연구: 음성 분류 솔루션"

# Row 32: update title (D32) and link (E32)
$ws.Range("D32").Value = "통계적 편향 (통계로 거짓말하기)"
$ws.Range("E32").Value = "https://dodonam.tistory.com/388"

# Row 51: update title (D51) and link (E51)
$ws.Range("D51").Value = "단순회귀분석 vs 다항회귀분석 vs 다중회귀분석"
$ws.Range("E51").Value = "https://bskyvision.com/entry/%EB%8B%A8%EC%88%9C%ED%9A%8C%EA%B7%80%EB%B6%84%EC%84%9D-vs-%EB%8B%A4%ED%95%AD%ED%9A%8C%EA%B7%80%EB%B6%84%EC%84%9D-vs-%EB%8B%A4%EC%A4%91%ED%9A%8C%EA%B7%80%EB%B6%84%EC%84%9D"
